$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.436.34'
$ws.Range("D3").Value = '1.850.75'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.27'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6292'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07644'
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("D11").Value = '2.072.25'
$ws.Range("E11").Value = '  +12.10%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.034'
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001063'
$ws.Range("E15").Value = '  -3.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.45'
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.185'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '29.521.44'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.51'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.33'
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.454'
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '157.58'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.436'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.70'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.386'
$ws.Range("E28").Value = '  +6.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.462'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05608'
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.053'
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.164'
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7026'
$ws.Range("E35").Value = '  -1.17%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = '1.231.02'
$ws.Range("E38").Value = '  -0.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.723'
$ws.Range("E39").Value = '  -1.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.445'
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9070'
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.98'
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.11'
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.198'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4026'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1157'
$ws.Range("E48").Value = '  +3.24%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.012'
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.682'
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05706'
$ws.Range("E51").Value = '  -0.02%  '
